$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.997.99'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.886.19'
$ws.Range("E3").Value = '  +1.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5152'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3749'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07185'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.09'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8993'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07648'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.28%  '
$ws.Range("D13").Value = '1.897.50'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.238'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008482'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '27.077.28'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.056'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").Value = '2.118.76'
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.375'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.296'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.730'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.912'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.783'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09203'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05040'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.230'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7645'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.993'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.279'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.594'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5626'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.069'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.100'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.629'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1501'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4816'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.79%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.601'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.54%  '
